# Alvearie FHIR IG - episode-duration StructureDefinition spreadsheet update
# Mirrors the gh-pages deploy: version bump 5.0.0 -> 6.0.0, refreshed date,
# Publisher/Jurisdiction metadata replacing the duplicated old "Contact" row,
# and the root Extension element's Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# -- Version (row 3) --
$meta.Range("B3").Value = "6.0.0"

# -- Date (row 8) --
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# -- Publisher value (row 9) was blank, now populated --
$meta.Range("B9").Value = "Alvearie Team"

# -- Row 10 used to duplicate row 11's "Contact" / "No display for
#    ContactDetail" pair. Delete it; everything below shifts up one row,
#    preserving cell types (text "false", etc.) instead of retyping them. --
$meta.Rows.Item(10).Delete()

# -- The row that is now 10 (old row 11, the other "Contact" copy) becomes
#    the new "Jurisdiction" / "United States of America" row. --
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# -- Elements sheet: root Extension element's Short/Definition text --
$elements.Range("K2").Value = "Episode Duration (Days)"
$elements.Range("L2").Value = "Duration of the episode of care, specified in days"
